# Auto-applies the cell-value edits described by the commit diff for cryptos.xlsx.
# (Updated symbol list / prices on Fri Dec 30 15:26:15 UTC 2022 with GitHub Actions)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'25.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.140"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05662"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'2.928"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8127"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8309"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1332"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.06956"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.02834"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.09390"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001516"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.0005956"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14OneONE"
$ws.Range("D16").Value = "'0.006107"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'3.501"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.3201"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.03164"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.1319"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04659"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.1358"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.001237"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004263"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009696"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001964"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03625"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1052"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002721"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003302"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.007376"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005288"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.2198"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "'0.002286"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
